# Rename the three logo pictures that live in the document's headers and
# footers (no visual/content change - just the drawing's display name,
# same as using Word's Selection Pane "Rename" on a picture):
#
#   Pearson logo (alt text ends in PearsonLogo.png), found once in each
#   footer story -> renamed from "image1.png" to "image2.png"
#
#   BTEC logo (alt text "BTec_Logo-Orange"), found in the "first page"
#   header story -> renamed from "image2.jpg" to "image1.jpg"
#
# Picture objects living in header/footer stories don't always honour a
# direct ".Name = ..." through the HeaderFooter.Range.InlineShapes path,
# so each shape is selected first and then renamed through
# $word.Selection.InlineShapes - mirroring how this rename is normally
# performed interactively (select the picture, then rename it).

$d = $word.ActiveDocument

for ($k = 1; $k -le $d.Sections.Count; $k++) {
    $sec = $d.Sections.Item($k)

    # --- Footers: Pearson logo, image1.png -> image2.png ---------------
    for ($i = 1; $i -le $sec.Footers.Count; $i++) {
        $ftr = $sec.Footers.Item($i)
        if ($ftr.Exists) {
            for ($j = 1; $j -le $ftr.Range.InlineShapes.Count; $j++) {
                $shp = $ftr.Range.InlineShapes.Item($j)
                if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                    [void]$shp.Select()
                    $word.Selection.InlineShapes.Item(1).Name = "image2.png"
                }
            }
        }
    }

    # --- Headers: BTEC logo, image2.jpg -> image1.jpg -------------------
    for ($i = 1; $i -le $sec.Headers.Count; $i++) {
        $hdr = $sec.Headers.Item($i)
        if ($hdr.Exists) {
            for ($j = 1; $j -le $hdr.Range.InlineShapes.Count; $j++) {
                $shp = $hdr.Range.InlineShapes.Item($j)
                if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    [void]$shp.Select()
                    $word.Selection.InlineShapes.Item(1).Name = "image1.jpg"
                }
            }
        }
    }
}
